# CIFAR_Result.xlsx — add final CIFAR-100 results (rows 28-33, cols H-O)
# and update the view/selection, per the commit "Add final results and
# visualiztion of CIFAR-100".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 28 (Resnet34): fill in the missing Fed/I,K,M,O columns ---
$ws.Range("I28").Value = 0.38233298562824503
$ws.Range("K28").Value = 0.38400000000000001
$ws.Range("M28").Value = 0.38601445988060701
$ws.Range("O28").Value = 0.38400000000000001

# --- Row 29 (Resnet50) ---
$ws.Range("H29").Value = 0.35498814025467101
$ws.Range("I29").Value = 0.35090711441950301
$ws.Range("J29").Value = 0.35980000000000001
$ws.Range("K29").Value = 0.35299999999999998
$ws.Range("L29").Value = 0.36361516869767802
$ws.Range("M29").Value = 0.35917234233016798
$ws.Range("N29").Value = 0.35980000000000001
$ws.Range("O29").Value = 0.35299999999999998

# --- Row 30 (VGG11) ---
$ws.Range("H30").Value = 0.43264137058639701
$ws.Range("I30").Value = 0.42309702160504398
$ws.Range("J30").Value = 0.43240000000000001
$ws.Range("K30").Value = 0.42809999999999998
$ws.Range("L30").Value = 0.45040682258230003
$ws.Range("M30").Value = 0.42833275179350699
$ws.Range("N30").Value = 0.43240000000000001
$ws.Range("O30").Value = 0.42809999999999998

# --- Row 31 (VGG13) ---
$ws.Range("H31").Value = 0.43179732631295198
$ws.Range("I31").Value = 0.399253818183811
$ws.Range("J31").Value = 0.43230000000000002
$ws.Range("K31").Value = 0.40510000000000002
$ws.Range("L31").Value = 0.44802092935657301
$ws.Range("M31").Value = 0.40223655995415403
$ws.Range("N31").Value = 0.43230000000000002
$ws.Range("O31").Value = 0.40510000000000002

# --- Row 32 (VGG16) ---
$ws.Range("H32").Value = 0.43692017285649998
$ws.Range("I32").Value = 0.36947713033653101
$ws.Range("J32").Value = 0.44040000000000001
$ws.Range("K32").Value = 0.37580000000000002
$ws.Range("L32").Value = 0.445718087135696
$ws.Range("M32").Value = 0.376053887413137
$ws.Range("N32").Value = 0.44040000000000001
$ws.Range("O32").Value = 0.37580000000000002

# --- Row 33 (VGG19) ---
$ws.Range("H33").Value = 0.45003995025529597
$ws.Range("I33").Value = 0.31733321055350699
$ws.Range("J33").Value = 0.44669999999999999
$ws.Range("K33").Value = 0.3387
$ws.Range("L33").Value = 0.463360608267949
$ws.Range("M33").Value = 0.325901286080852
$ws.Range("N33").Value = 0.44669999999999999
$ws.Range("O33").Value = 0.3387

# --- View: scroll/select like the saved author session ---
$ws.Range("F41").Select()
